# Auto-generated edit script applying the diff to before.xlsx
$wb = $excel.ActiveWorkbook

# ---- Sheet1 ("iModulon"): update column B values for rows 2-14 ----
$ws1 = $wb.Worksheets.Item("iModulon")
$ws1.Range("B2").Value = 1.104972375690608
$ws1.Range("B3").Value = 1.657458563535912
$ws1.Range("B4").Value = 8.103130755064457
$ws1.Range("B5").Value = 7.366482504604052
$ws1.Range("B6").Value = 0.1841620626151013
$ws1.Range("B7").Value = 2.762430939226519
$ws1.Range("B8").Value = 1.841620626151013
$ws1.Range("B9").Value = 0.9208103130755065
$ws1.Range("B10").Value = 16.02209944751381
$ws1.Range("B11").Value = 0.7366482504604052
$ws1.Range("B12").Value = 0.5524861878453038
$ws1.Range("B13").Value = 0.9208103130755065
$ws1.Range("B14").Value = 3.683241252302026

# ---- Sheet2 ("Subsystem"): rewrite full table (52 -> 56 data rows) ----
$ws2 = $wb.Worksheets.Item("Subsystem")

# Pre-format the 4 newly-needed rows (53-56) by copying the category-name
# column style (bold/border/centered) from an existing formatted cell.
$ws2.Range("A2").Copy()
$ws2.Range("A53:A56").PasteSpecial(-4122)

$ws2.Range("A2").Value = "Alanine, aspartate and glutamate metabolism"
$ws2.Range("B2").Value = 1.557093425605536
$ws2.Range("C2").Value = 2.083333333333333
$ws2.Range("A3").Value = "Aminosugars metabolism"
$ws2.Range("B3").Value = 0.8650519031141869
$ws2.Range("C3").Value = 0.6944444444444444
$ws2.Range("A4").Value = "Arginine and proline metabolism"
$ws2.Range("B4").Value = 2.768166089965398
$ws2.Range("C4").Value = 9.722222222222223
$ws2.Range("A5").Value = "Biomass and maintenance functions"
$ws2.Range("B5").Value = 0.1730103806228374
$ws2.Range("C5").Value = 0
$ws2.Range("A6").Value = "Biotin metabolism"
$ws2.Range("B6").Value = 0
$ws2.Range("C6").Value = 2.083333333333333
$ws2.Range("A7").Value = "C5-Branched dibasic acid metabolism"
$ws2.Range("B7").Value = 0.5190311418685121
$ws2.Range("C7").Value = 0
$ws2.Range("A8").Value = "Carbon fixation"
$ws2.Range("B8").Value = 1.038062283737024
$ws2.Range("C8").Value = 0
$ws2.Range("A9").Value = "Carotenoid Biosynthesis"
$ws2.Range("B9").Value = 1.903114186851211
$ws2.Range("C9").Value = 0
$ws2.Range("A10").Value = "Citrate cycle (TCA cycle)"
$ws2.Range("B10").Value = 1.038062283737024
$ws2.Range("C10").Value = 0
$ws2.Range("A11").Value = "Cyanophycin metabolism"
$ws2.Range("B11").Value = 0.3460207612456748
$ws2.Range("C11").Value = 0
$ws2.Range("A12").Value = "Extracellular exchange"
$ws2.Range("B12").Value = 3.28719723183391
$ws2.Range("C12").Value = 0
$ws2.Range("A13").Value = "Fatty acid biosynthesis"
$ws2.Range("B13").Value = 17.47404844290658
$ws2.Range("C13").Value = 0
$ws2.Range("A14").Value = "Folate biosynthesis"
$ws2.Range("B14").Value = 2.249134948096886
$ws2.Range("C14").Value = 1.388888888888889
$ws2.Range("A15").Value = "Fructose and mannose metabolism"
$ws2.Range("B15").Value = 0
$ws2.Range("C15").Value = 5.555555555555555
$ws2.Range("A16").Value = "Galactolipids metabolism"
$ws2.Range("B16").Value = 4.671280276816609
$ws2.Range("C16").Value = 0
$ws2.Range("A17").Value = "Glutamate metabolism"
$ws2.Range("B17").Value = 1.211072664359861
$ws2.Range("C17").Value = 0.6944444444444444
$ws2.Range("A18").Value = "Glutathione metabolism"
$ws2.Range("B18").Value = 0.1730103806228374
$ws2.Range("C18").Value = 1.388888888888889
$ws2.Range("A19").Value = "Glycerolipid metabolism"
$ws2.Range("B19").Value = 0.1730103806228374
$ws2.Range("C19").Value = 1.388888888888889
$ws2.Range("A20").Value = "Glycolysis/Gluconeogenesis"
$ws2.Range("B20").Value = 2.941176470588235
$ws2.Range("C20").Value = 3.472222222222222
$ws2.Range("A21").Value = "Glyoxylate and dicarboxylate metabolism"
$ws2.Range("B21").Value = 2.076124567474048
$ws2.Range("C21").Value = 1.388888888888889
$ws2.Range("A22").Value = "Histidine metabolism"
$ws2.Range("B22").Value = 1.557093425605536
$ws2.Range("C22").Value = 1.388888888888889
$ws2.Range("A23").Value = "Hydrogen production"
$ws2.Range("B23").Value = 0
$ws2.Range("C23").Value = 0.6944444444444444
$ws2.Range("A24").Value = "Inositol phosphate metabolism"
$ws2.Range("B24").Value = 0
$ws2.Range("C24").Value = 1.388888888888889
$ws2.Range("A25").Value = "Intracellular demand"
$ws2.Range("B25").Value = 0.1730103806228374
$ws2.Range("C25").Value = 0
$ws2.Range("A26").Value = "Intracellular source/sink"
$ws2.Range("B26").Value = 0.1730103806228374
$ws2.Range("C26").Value = 0
$ws2.Range("A27").Value = "Lipopolysaccharide biosynthesis"
$ws2.Range("B27").Value = 0.8650519031141869
$ws2.Range("C27").Value = 0
$ws2.Range("A28").Value = "Lysine metabolism"
$ws2.Range("B28").Value = 1.557093425605536
$ws2.Range("C28").Value = 1.388888888888889
$ws2.Range("A29").Value = "Nicotinate and nicotinamide metabolism"
$ws2.Range("B29").Value = 0.8650519031141869
$ws2.Range("C29").Value = 0.6944444444444444
$ws2.Range("A30").Value = "Nitrogen metabolism"
$ws2.Range("B30").Value = 1.730103806228374
$ws2.Range("C30").Value = 0.6944444444444444
$ws2.Range("A31").Value = "Nucleotide sugars metabolism"
$ws2.Range("B31").Value = 0
$ws2.Range("C31").Value = 2.777777777777778
$ws2.Range("A32").Value = "Others"
$ws2.Range("B32").Value = 0.6920415224913495
$ws2.Range("C32").Value = 1.388888888888889
$ws2.Range("A33").Value = "Oxidative phosphorylation"
$ws2.Range("B33").Value = 1.211072664359861
$ws2.Range("C33").Value = 2.083333333333333
$ws2.Range("A34").Value = "PHB byosynthesis"
$ws2.Range("B34").Value = 0
$ws2.Range("C34").Value = 2.083333333333333
$ws2.Range("A35").Value = "Pantothenate and CoA biosynthesis"
$ws2.Range("B35").Value = 1.557093425605536
$ws2.Range("C35").Value = 0
$ws2.Range("A36").Value = "Pentose phosphate pathway"
$ws2.Range("B36").Value = 1.384083044982699
$ws2.Range("C36").Value = 0.6944444444444444
$ws2.Range("A37").Value = "Peptidoglycan biosynthesis"
$ws2.Range("B37").Value = 1.384083044982699
$ws2.Range("C37").Value = 0.6944444444444444
$ws2.Range("A38").Value = "Phenylalanine tyrosine and tryptophan biosynthesis"
$ws2.Range("B38").Value = 3.28719723183391
$ws2.Range("C38").Value = 5.555555555555555
$ws2.Range("A39").Value = "Photosynthesis"
$ws2.Range("B39").Value = 1.038062283737024
$ws2.Range("C39").Value = 0
$ws2.Range("A40").Value = "Porphyrin and chlorophyll metabolism"
$ws2.Range("B40").Value = 7.093425605536333
$ws2.Range("C40").Value = 4.861111111111112
$ws2.Range("A41").Value = "Purine metabolism"
$ws2.Range("B41").Value = 4.325259515570934
$ws2.Range("C41").Value = 4.861111111111112
$ws2.Range("A42").Value = "Pyrimidine metabolism"
$ws2.Range("B42").Value = 3.28719723183391
$ws2.Range("C42").Value = 2.083333333333333
$ws2.Range("A43").Value = "Pyruvate metabolism"
$ws2.Range("B43").Value = 1.730103806228374
$ws2.Range("C43").Value = 1.388888888888889
$ws2.Range("A44").Value = "Riboflavin metabolism"
$ws2.Range("B44").Value = 1.730103806228374
$ws2.Range("C44").Value = 0
$ws2.Range("A45").Value = "Starch and sucrose metabolism"
$ws2.Range("B45").Value = 0.8650519031141869
$ws2.Range("C45").Value = 3.472222222222222
$ws2.Range("A46").Value = "Steroid biosynthesis"
$ws2.Range("B46").Value = 0
$ws2.Range("C46").Value = 2.083333333333333
$ws2.Range("A47").Value = "Sterol biosynthesis"
$ws2.Range("B47").Value = 1.038062283737024
$ws2.Range("C47").Value = 0
$ws2.Range("A48").Value = "Sulfolipid Biosynthesis"
$ws2.Range("B48").Value = 1.730103806228374
$ws2.Range("C48").Value = 0
$ws2.Range("A49").Value = "Sulfur Cysteine and methionine metabolism"
$ws2.Range("B49").Value = 2.422145328719723
$ws2.Range("C49").Value = 3.472222222222222
$ws2.Range("A50").Value = "Terpenoid backbone biosynthesis"
$ws2.Range("B50").Value = 1.557093425605536
$ws2.Range("C50").Value = 0
$ws2.Range("A51").Value = "Thiamine metabolism"
$ws2.Range("B51").Value = 0
$ws2.Range("C51").Value = 3.472222222222222
$ws2.Range("A52").Value = "Transport"
$ws2.Range("B52").Value = 6.920415224913495
$ws2.Range("C52").Value = 15.97222222222222
$ws2.Range("A53").Value = "Ubiquinone and other pterpenoids biosynthesis"
$ws2.Range("B53").Value = 2.422145328719723
$ws2.Range("C53").Value = 0
$ws2.Range("A54").Value = "Urea cycle and metabolism of amino groups"
$ws2.Range("B54").Value = 1.211072664359861
$ws2.Range("C54").Value = 0.6944444444444444
$ws2.Range("A55").Value = "Valine leucine and isoleucine biosynthesis"
$ws2.Range("B55").Value = 1.730103806228374
$ws2.Range("C55").Value = 2.777777777777778
$ws2.Range("A56").Value = "Vitamin B6 metabolism"
$ws2.Range("B56").Value = 0
$ws2.Range("C56").Value = 3.472222222222222
